$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Isle of Man "Investments & Asset Management" comps table was refreshed:
#  - the PME African Infrastructure Opportunities row (old row 5) was dropped
#  - the Origo Partners row (old row 6) moved up to row 5 with refreshed figures
#  - every other company row keeps its position but gets refreshed metrics

# Row 2 -- updated metrics (company id "4" -> "3")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("D2").Value = -0.293
$ws.Range("G2").Value = 0.4237879316358563
$ws.Range("H2").Value = 0.4237879316358563
$ws.Range("I2").Value = -0.2040922322557597
$ws.Range("J2").Value = -0.2040922322557597
$ws.Range("K2").Value = 0.705
$ws.Range("L2").Value = 0.2459016393442623
$ws.Range("M2").Value = 2.1
$ws.Range("N2").Value = 0.0248370805785857
$ws.Range("O2").Value = 2.978723404255319
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 0.0248370805785857
$ws.Range("R2").Value = 2.978723404255319
$ws.Range("U2").Value = 7.67
$ws.Range("V2").Value = 0.09071448001797731
$ws.Range("W2").Value = 0.003814713896457766
$ws.Range("X2").Value = 0.03846338425444983
$ws.Range("Y2").Value = -0.03464867035799206
$ws.Range("Z2").Value = 0.4484128083425841
$ws.Range("AA2").Value = -0.05266622778143517
$ws.Range("AB2").Value = 0.03846338425444983
$ws.Range("AC2").Value = -0.09112961203588499
$ws.Range("AD2").Value = 0.546
$ws.Range("AE2").Value = 0.06566214938631565
$ws.Range("AF2").Value = 0.6116621493863157
$ws.Range("AG2").Value = -7.058337850613684
$ws.Range("AH2").Value = 0.007182280754838091
$ws.Range("AI2").Value = 0.01727778054050829
$ws.Range("AJ2").Value = -0.09108395111019657
$ws.Range("AK2").Value = -0.2545227117145548
$ws.Range("AL2").Value = 0.062
$ws.Range("AM2").Value = -0.014
$ws.Range("AN2").Value = -1.649546827794562
$ws.Range("AO2").Value = -10.25806451612903
$ws.Range("AP2").Value = 21.32428353659724
$ws.Range("AQ2").Value = 45.42857142857143

# Row 3 -- Cambria Africa plc (AIM:CMB): updated metrics
$ws.Range("D3").Value = -0.293
$ws.Range("G3").Value = 1.098958333333333
$ws.Range("H3").Value = 1.098958333333333
$ws.Range("I3").Value = 0.03638935943892546
$ws.Range("J3").Value = 0.03638935943892546
$ws.Range("K3").Value = 0.028
$ws.Range("L3").Value = 0.01458333333333333
$ws.Range("U3").Value = 2.33
$ws.Range("V3").Value = 1.226315789473684
$ws.Range("W3").Value = 0.003814713896457766
$ws.Range("X3").Value = 0.04785187809468971
$ws.Range("Y3").Value = -0.04403716419823195
$ws.Range("Z3").Value = 0.4429410952527929
$ws.Range("AA3").Value = 0.0161183427254252
$ws.Range("AB3").Value = 0.04378244189685672
$ws.Range("AC3").Value = -0.02766409917143152
$ws.Range("AD3").Value = 0.546
$ws.Range("AE3").Value = 0.06566214938631565
$ws.Range("AF3").Value = 0.6116621493863157
$ws.Range("AG3").Value = -1.718337850613684
$ws.Range("AH3").Value = 0.2435288319074942
$ws.Range("AI3").Value = 0.07376833961229187
$ws.Range("AJ3").Value = -9.458975666744617
$ws.Range("AK3").Value = -0.2882313367574121
$ws.Range("AL3").Value = 0.062
$ws.Range("AM3").Value = 0.062
$ws.Range("AN3").Value = 2.256198347107438
$ws.Range("AO3").Value = 0.3064516129032258
$ws.Range("AP3").Value = -7.10056963063506
$ws.Range("AQ3").Value = 0.3064516129032258

# Row 4 -- Agronomics Limited (AIM:ANIC): updated metrics
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = -0.09852216748768472
$ws.Range("J4").Value = -0.09852216748768472
$ws.Range("K4").Value = 0.758
$ws.Range("L4").Value = 0.9334975369458127
$ws.Range("O4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 3.45
$ws.Range("V4").Value = 0.04217603911980441
$ws.Range("W4").Value = 0.3697560975609757
$ws.Range("X4").Value = 0.03846338425444983
$ws.Range("Y4").Value = 0.3312927133065258
$ws.Range("Z4").Value = 0.5345622119815669
$ws.Range("AA4").Value = -0.05266622778143517
$ws.Range("AB4").Value = 0.03846338425444983
$ws.Range("AC4").Value = -0.09112961203588499
$ws.Range("AG4").Value = -3.45
$ws.Range("AJ4").Value = -0.04403318442884493
$ws.Range("AK4").Value = -0.1678832116788321
$ws.Range("AM4").Value = -0.076
$ws.Range("AQ4").Value = 1.052631578947368

# Row 5 -- now Origo Partners PLC (AIM:OPP), replacing PME African Infrastructure Opportunities PLC
$ws.Range("B5").Value = "Origo Partners PLC (AIM:OPP)"
$ws.Range("G5").Value = -6.62962962962963
$ws.Range("H5").Value = -6.62962962962963
$ws.Range("I5").Value = -4.259259259259259
$ws.Range("J5").Value = -4.259259259259259
$ws.Range("K5").Value = -0.081
$ws.Range("L5").Value = -0.6
$ws.Range("M5").Value = 2.1
$ws.Range("N5").Value = 2.467685076380729
$ws.Range("O5").Value = -25.92592592592593
$ws.Range("P5").Value = 2.1
$ws.Range("Q5").Value = 2.467685076380729
$ws.Range("R5").Value = -25.92592592592593
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1.89
$ws.Range("V5").Value = 2.220916568742656
$ws.Range("W5").Value = -0.01519699812382739
$ws.Range("X5").Value = 0.03846338425444983
$ws.Range("Y5").Value = -0.05366038237827722
$ws.Range("Z5").Value = 0.25
$ws.Range("AA5").Value = -1.064814814814815
$ws.Range("AB5").Value = 0.03846338425444983
$ws.Range("AC5").Value = -1.103278199069264
$ws.Range("AD5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -1.89
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 1.819056785370549
$ws.Range("AK5").Value = -1.549180327868853
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = -0
$ws.Range("AP5").Value = 3.298429319371728
$ws.Range("AO5").ClearContents()
$ws.Range("AQ5").ClearContents()

# Drop the old row 6 (its data now lives in row 5 above)
$ws.Rows(6).Delete()
